# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(5, 1, 2, 5, 0, 2, 2, 1, 3, 2, 3, 4, 4, 3, 7, 7, 3, 6, 4, 7, 9, 3, 2, 3, 4)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $kValues[$i]
}
